# Apply the cryptos-list refresh described in the commit: updates the
# Price (D) / Volume(1h) (E) columns for most rows, and for rows 10-11
# the two coins (TRON / WrappedEther) swap places along with their data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "25.491.19", "0.9988") that must
# stay text, exactly as stored in the workbook. A leading apostrophe forces
# Excel to keep it as text instead of auto-coercing it to a number.

$ws.Range('D2').Value = '''25.491.19'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '''1.668.93'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('D4').Value = '''0.9988'
$ws.Range('D5').Value = '''238.23'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '''0.4801'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('D9').Value = '''0.06174'
$ws.Range('E9').Value = '  +2.87%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').Value = '''0.06990'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '''1.666.35'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '''0.5895'
$ws.Range('E13').Value = '  -4.81%  '
$ws.Range('E14').Value = '  -3.05%  '
$ws.Range('D15').Value = '''75.02'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D17').Value = '''0.9996'
$ws.Range('D18').Value = '''25.484.63'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = '''0.000006766'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = '''1.881.23'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('D22').Value = '''4.455'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = '''8.739'
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('D24').Value = '''5.283'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '''136.97'
$ws.Range('E25').Value = '  +3.70%  '
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('D27').Value = '''1.392'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '''1.726'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').Value = '''104.92'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('E30').Value = '  +5.61%  '
$ws.Range('D31').Value = '''0.07804'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').Value = '''3.652'
$ws.Range('D33').Value = '''0.9989'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '''0.04251'
$ws.Range('E34').Value = '  -4.76%  '
$ws.Range('D35').Value = '''2.599'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('E36').Value = '  +4.36%  '
$ws.Range('D37').Value = '''0.9496'
$ws.Range('E37').Value = '  +2.20%  '
$ws.Range('D38').Value = '''2.591'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '''0.8585'
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  -5.49%  '
$ws.Range('D42').Value = '''1.852'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = '''96.16'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').Value = '''0.3773'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = '''4.827'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D46').Value = '''0.1118'
$ws.Range('E46').Value = '  -2.65%  '
$ws.Range('D47').Value = '''6.197'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').Value = '''0.05248'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '''29.83'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').Value = '''7.376'
$ws.Range('E50').Value = '  +3.11%  '
$ws.Range('E51').Value = '  +0.19%  '
